$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to Text format so numeric-looking strings
# (e.g. "323.40", "0.0790", "1.90") keep their exact original formatting
# instead of being auto-converted into floating point numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '48.262.58'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '2.518.66'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '323.40'
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = '109.26'
$ws.Range("E6").Value = '  -0.57%  '
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.561'
$ws.Range("E9").Value = '  +4.25%  '
$ws.Range("D10").Value = '40.30'
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("D11").Value = '19.87'
$ws.Range("E11").Value = '  +6.93%  '
$ws.Range("E13").Value = '  +0.80%  '
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("D15").Value = '2.913.23'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("D16").Value = '2.524.14'
$ws.Range("E16").Value = '  +0.66%  '
$ws.Range("E17").Value = '  -0.70%  '
$ws.Range("D18").Value = '48.154.11'
$ws.Range("E18").Value = '  +1.72%  '
$ws.Range("D19").Value = '13.36'
$ws.Range("E19").Value = '  +3.81%  '
$ws.Range("E20").Value = '  -0.70%  '
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  +1.64%  '
$ws.Range("D23").Value = '72.64'
$ws.Range("E23").Value = '  +2.88%  '
$ws.Range("D24").Value = '269.74'
$ws.Range("E24").Value = '  +8.22%  '
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").Value = '26.17'
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '0.145'
$ws.Range("E29").Value = '  +4.65%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '2.20'
$ws.Range("E30").Value = '  -4.06%  '
$ws.Range("D31").Value = '35.21'
$ws.Range("D32").Value = '49.81'
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = '19.99'
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("D34").Value = '5.41'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '0.0790'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("D38").Value = '4.72'
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '22.22'
$ws.Range("E41").Value = '  +4.33%  '
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").Value = '118.44'
$ws.Range("E43").Value = '  -2.89%  '
$ws.Range("D44").Value = '0.0300'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").Value = '2.000.44'
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = '1.90'
$ws.Range("E47").Value = '  +6.54%  '
$ws.Range("E48").Value = '  -2.27%  '
$ws.Range("E49").Value = '  +0.27%  '
$ws.Range("D50").Value = '5.26'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("D51").Value = '80.99'
$ws.Range("E51").Value = '  +3.28%  '
